$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.071.64"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "3.221.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'604.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'154.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "3.218.78"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.534"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'6.15"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.512"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'39.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.41%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "3.747.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'7.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.69%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "66.149.63"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "3.230.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'511.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'15.49"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.741"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("E22").ClearFormats()
$ws.Range("B23").Value = "Uniswap"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C23").ClearFormats()
$ws.Range("D23").Value = "'8.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.36%  "
$ws.Range("E23").ClearFormats()
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C24").ClearFormats()
$ws.Range("D24").Value = "'15.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'85.27"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'3.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'9.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'2.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("E30").ClearFormats()
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("B31").ClearFormats()
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C31").ClearFormats()
$ws.Range("D31").Value = "'6.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.30%  "
$ws.Range("E31").ClearFormats()
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = "'28.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'6.63"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'55.30"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.0907"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'484.38"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.0421"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E39").ClearFormats()
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C40").ClearFormats()
$ws.Range("D40").Value = "'2.97"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("E40").ClearFormats()
$ws.Range("B41").Value = "Cosmos"
$ws.Range("B41").ClearFormats()
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C41").ClearFormats()
$ws.Range("D41").Value = "'8.96"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.302"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.99%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.120"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E43").ClearFormats()
$ws.Range("B44").Value = "PEPE"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "0.0₃0649"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.44%  "
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = "Maker"
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "2.950.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'2.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'28.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'2.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'121.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.43%  "
$ws.Range("E51").ClearFormats()
